$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price values (price list update ~ +21%), column D, rows 2-32.
$rng = $ws.Range("D2:D32")

$ws.Range("D2").Value = 38.83
$ws.Range("D3").Value = 37.43
$ws.Range("D4").Value = 40.97
$ws.Range("D5").Value = 33.119999999999997
$ws.Range("D6").Value = 36.93
$ws.Range("D7").Value = 47.27
$ws.Range("D8").Value = 16.7
$ws.Range("D9").Value = 17.2
$ws.Range("D10").Value = 26.810000000000002
$ws.Range("D11").Value = 16.96
$ws.Range("D12").Value = 17.27
$ws.Range("D13").Value = 25.37
$ws.Range("D14").Value = 15.4
$ws.Range("D15").Value = 16.8
$ws.Range("D16").Value = 29.7
$ws.Range("D17").Value = 13.5
$ws.Range("D18").Value = 18.48
$ws.Range("D19").Value = 29.48
$ws.Range("D20").Value = 46.54
$ws.Range("D21").Value = 28.400000000000002
$ws.Range("D22").Value = 31.310000000000002
$ws.Range("D23").Value = 30.43
$ws.Range("D24").Value = 24.45
$ws.Range("D25").Value = 33.18
$ws.Range("D26").Value = 48.4
$ws.Range("D27").Value = 32.160000000000004
$ws.Range("D28").Value = 30.150000000000002
$ws.Range("D29").Value = 60.660000000000004
$ws.Range("D30").Value = 40.630000000000003
$ws.Range("D31").Value = 48.86
$ws.Range("D32").Value = 60.660000000000004

# New column D style: number format #,##0.00 (builtin 4), centered
# horizontally, top-aligned vertically (was center/center before).
$rng.NumberFormat = "#,##0.00"
$rng.VerticalAlignment = -4160
$rng.HorizontalAlignment = -4108

# Match the saved selection state (D2:D32 active, anchor D2).
$rng.Select()
